# Avg Vehicle Loading.xlsx - "Rail and aviation edits"
#
# Adds a new assumption row ("Assumption - train cars per locomotive" = 10)
# to the Freight Rail block on the "BTS NTS Modal Profile Data" sheet, and
# folds that multiplier into the weighted-average distance-weighted
# passengers/Amtrak-train formula. Inserting the row naturally re-points
# every downstream formula (same sheet and the two sheets that pull values
# from it) at their new row numbers.

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("BTS NTS Modal Profile Data")
$passSheet = $wb.Worksheets.Item("AVLo-passengers")
$freightSheet = $wb.Worksheets.Item("AVLo-freight")
$aboutSheet = $wb.Worksheets.Item("About")

# --- Insert the new assumption row above the old row 22 ("Freight Rail" block) ---
$dataSheet.Rows("22:22").Insert()

$dataSheet.Range("A22").Value = "Assumption - train cars per locomotive"
$dataSheet.Range("B22").Value = 10

# Match the plain "bold label" look used elsewhere (style index 1 in the
# original file) rather than the shaded header style that Insert copied
# down from the row above.
$dataSheet.Range("A22:B22").Font.Bold = $true
$dataSheet.Range("A22:B22").Interior.Pattern = -4142   # xlNone
$dataSheet.Range("A22:B22").Interior.ColorIndex = -4142 # xlColorIndexNone

# --- Fold the new assumption into the weighted-average formula (now row 37) ---
$dataSheet.Range("B37").Formula = "=(B26*B25+B34*B28+B35*B29+B36*B30)/SUM(B25,B28:B30)*B22"

# --- Update view state (active cell / selection) on each sheet ---
$aboutSheet.Activate()
$aboutSheet.Range("B11").Select()

$dataSheet.Activate()
$dataSheet.Range("A20:XFD20").Select()

$passSheet.Activate()
$passSheet.Rows.Item(1).RowHeight = 45
$passSheet.Range("B5").Select()

$freightSheet.Activate()
$freightSheet.Rows.Item(1).RowHeight = 45
$freightSheet.Range("B6").Select()

# Leave "About" as the active sheet/tab, matching the saved workbook state.
$aboutSheet.Activate()
